$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 20
$ws.Range("E8").Value = 28
$ws.Range("E12").Value = 18
$ws.Range("E14").Value = 29
$ws.Range("E15").Value = 69
$ws.Range("E16").Value = 233
$ws.Range("E18").Value = 66
